# daily auto push: 2026-01-21 13:53 UTC
# Insert a new data row at row 698 (pushing the existing rows 698-739 down
# to 699-740) and populate it with the new day's first sample:
#   2026/01/21, 水, 19, 201
# This mirrors the diff where dimension grows from A1:D739 to A1:D740 and
# all rows from 698 onward shift down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 698..739 down to 699..740, creating a fresh blank row 698.
$ws.Rows(698).Insert()

# Column A holds dates stored as literal text (e.g. "2026/12/29"), not as
# real Excel dates, elsewhere in the sheet. Force the same text storage
# here so "2026/01/21" isn't auto-converted into a date serial number.
$ws.Cells.Item(698, 1).NumberFormat = "@"
$ws.Cells.Item(698, 1).Value = "2026/01/21"
# Restore the plain/default style (no explicit number format) so the new
# cell doesn't end up with a stray text-format style applied to it.
$ws.Cells.Item(698, 1).Style = "Normal"

$ws.Cells.Item(698, 2).Value = "水"
$ws.Cells.Item(698, 3).Value = 19
$ws.Cells.Item(698, 4).Value = 201
